$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (old row 11) so everything below shifts up by one
# (22 data rows -> 21 data rows).
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: updated timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$ws.Range("B9").Value = "Alvearie Team"

# The former second "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value now populated with the literal text "true"
# (a bare Value = "true" is auto-typed as a boolean, so force text entry with a
# leading quote, then re-copy the surrounding cell's formatting over it so the
# cell keeps its normal, non "quote-prefixed" style).
$ws.Range("B14").Value = "'true"
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
